$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the header timestamp (also updates the Table column header)
$ws.Range("C1").Value = "2024-11-15 17:18:48"

# Update attendance status for specific rows: Falta/Retardo -> Puntual
$ws.Range("C3").Value = "Puntual"
$ws.Range("C5").Value = "Puntual"
$ws.Range("C6").Value = "Puntual"
$ws.Range("C7").Value = "Puntual"
$ws.Range("C8").Value = "Puntual"
$ws.Range("C10").Value = "Puntual"
$ws.Range("C12").Value = "Puntual"
